# chartink_screener "10per change" sheet — break out stock.yaml completed
# 1) D2:D10 (bsecode) were entered as text; they should be plain numbers.
# 2) Rows 11-27 are new stock rows appended below the existing table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Fix D2:D10 so the BSE code is numeric, not text ---------------
$bsecodes = @{
    2  = 500408
    3  = 533179
    4  = 500114
    5  = 543237
    6  = 533758
    7  = 500770
    8  = 542830
    9  = 500271
    10 = 540133
}
foreach ($r in $bsecodes.Keys) {
    $ws.Cells.Item($r, 4).Value = $bsecodes[$r]
}

# --- 2. Append the new rows (11-27) ------------------------------------
$newRows = @(
    @(1,  "DIXON",      "Dixon Technologies",                     "540699", -0.6899999999999999,  8790,    64587,   "05/06/2024 04:07:17"),
    @(2,  "LT",         "Larsen & Toubro Limited",                "500510", -3.52,                 3283.35, 1212471, "05/06/2024 04:07:17"),
    @(3,  "GRASIM",     "Grasim Industries Limited",              "500300", -1.17,                 2229.9,  105697,  "05/06/2024 04:07:17"),
    @(4,  "HAVELLS",    "Havells India Limited",                  "517354", -0.3,                  1767.6,  99780,   "05/06/2024 04:07:17"),
    @(5,  "PRESTIGE",   "Prestige Estates Projects Limited",      "533274", -2.65,                 1471.75, 124151,  "05/06/2024 04:07:17"),
    @(6,  "BHARATFORG", "Bharat Forge Limited",                   "500493", -0.43,                 1440,    267731,  "05/06/2024 04:07:17"),
    @(7,  "VBL",        "Varun Beverages Limited",                "540180", -1.11,                 1397.7,  221982,  "05/06/2024 04:07:17"),
    @(8,  "ATGL",       "Adani Total Gas Ltd",                    "542066", -1.76,                 892.7,   1736794, "05/06/2024 04:07:17"),
    @(9,  "SBIN",       "State Bank Of India",                    "500112", -0.73,                 769.55,  8741287, "05/06/2024 04:07:17"),
    @(10, "ADANIPOWER", "Adani Power Limited",                    "533096", -3.69,                 696.3,   5914928, "05/06/2024 04:07:17"),
    @(11, "HINDALCO",   "Hindalco Industries Limited",            "500440", -4.42,                 620.35,  5061288, "05/06/2024 04:07:17"),
    @(12, "CGPOWER",    "CG Power and Industrial Solutions Ltd",  "500093", -4.39,                 599.45,  1195141, "05/06/2024 04:07:17"),
    @(13, "FORTIS",     "Fortis Healthcare Limited",              "532843", -1.42,                 433.1,   56348,   "05/06/2024 04:07:17"),
    @(14, "PAYTM",      "One 97 Communications Ltd",              "543396", -4.99,                 339.55,  2424683, "05/06/2024 04:07:17"),
    @(15, "POWERGRID",  "Power Grid Corporation Of India Limited","532898", -2.1,                  289.75,  7538472, "05/06/2024 04:07:17"),
    @(16, "ASHOKLEY",   "Ashok Leyland Limited",                  "500477", 1.42,                  210.7,   3448763, "05/06/2024 04:07:17"),
    @(17, "GAIL",       "Gail (india) Limited",                   "532155", -0.63,                 189.1,   8946102, "05/06/2024 04:07:17")
)

$rowIndex = 11
foreach ($row in $newRows) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]

    # bsecode stays text (as in the source data) — force text so Excel
    # doesn't coerce the numeric-looking string into a number.
    $dCell = $ws.Cells.Item($rowIndex, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[3]

    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws.Cells.Item($rowIndex, 8).Value = $row[7]

    $rowIndex++
}
